# Retraining the 15 minutes for Elnet
# Updates the forecasted consumption (A) and timestamp (B) columns: data refreshed
# and timestamps shifted forward by 4 days (45919 -> 45923).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 5160
$ws.Cells.Item(2, 2).Value = 45923
$ws.Cells.Item(3, 1).Value = 5110
$ws.Cells.Item(3, 2).Value = 45923.01041666666
$ws.Cells.Item(4, 1).Value = 5080
$ws.Cells.Item(4, 2).Value = 45923.02083333334
$ws.Cells.Item(5, 1).Value = 5060
$ws.Cells.Item(5, 2).Value = 45923.03125
$ws.Cells.Item(6, 1).Value = 5040
$ws.Cells.Item(6, 2).Value = 45923.04166666666
$ws.Cells.Item(7, 1).Value = 5030
$ws.Cells.Item(7, 2).Value = 45923.05208333334
$ws.Cells.Item(8, 1).Value = 5020
$ws.Cells.Item(8, 2).Value = 45923.0625
$ws.Cells.Item(9, 1).Value = 5020
$ws.Cells.Item(9, 2).Value = 45923.07291666666
$ws.Cells.Item(10, 1).Value = 5000
$ws.Cells.Item(10, 2).Value = 45923.08333333334
$ws.Cells.Item(11, 1).Value = 5000
$ws.Cells.Item(11, 2).Value = 45923.09375
$ws.Cells.Item(12, 1).Value = 5000
$ws.Cells.Item(12, 2).Value = 45923.10416666666
$ws.Cells.Item(13, 1).Value = 5000
$ws.Cells.Item(13, 2).Value = 45923.11458333334
$ws.Cells.Item(14, 1).Value = 5000
$ws.Cells.Item(14, 2).Value = 45923.125
$ws.Cells.Item(15, 1).Value = 5020
$ws.Cells.Item(15, 2).Value = 45923.13541666666
$ws.Cells.Item(16, 1).Value = 5030
$ws.Cells.Item(16, 2).Value = 45923.14583333334
$ws.Cells.Item(17, 1).Value = 5060
$ws.Cells.Item(17, 2).Value = 45923.15625
$ws.Cells.Item(18, 1).Value = 5110
$ws.Cells.Item(18, 2).Value = 45923.16666666666
$ws.Cells.Item(19, 1).Value = 5170
$ws.Cells.Item(19, 2).Value = 45923.17708333334
$ws.Cells.Item(20, 1).Value = 5250
$ws.Cells.Item(20, 2).Value = 45923.1875
$ws.Cells.Item(21, 1).Value = 5350
$ws.Cells.Item(21, 2).Value = 45923.19791666666
$ws.Cells.Item(22, 1).Value = 5460
$ws.Cells.Item(22, 2).Value = 45923.20833333334
$ws.Cells.Item(23, 1).Value = 5590
$ws.Cells.Item(23, 2).Value = 45923.21875
$ws.Cells.Item(24, 1).Value = 5710
$ws.Cells.Item(24, 2).Value = 45923.22916666666
$ws.Cells.Item(25, 1).Value = 5860
$ws.Cells.Item(25, 2).Value = 45923.23958333334
$ws.Cells.Item(26, 1).Value = 6060
$ws.Cells.Item(26, 2).Value = 45923.25
$ws.Cells.Item(27, 1).Value = 6180
$ws.Cells.Item(27, 2).Value = 45923.26041666666
$ws.Cells.Item(28, 1).Value = 6270
$ws.Cells.Item(28, 2).Value = 45923.27083333334
$ws.Cells.Item(29, 1).Value = 6300
$ws.Cells.Item(29, 2).Value = 45923.28125
$ws.Cells.Item(30, 1).Value = 6300
$ws.Cells.Item(30, 2).Value = 45923.29166666666
$ws.Cells.Item(31, 1).Value = 6290
$ws.Cells.Item(31, 2).Value = 45923.30208333334
$ws.Cells.Item(32, 1).Value = 6250
$ws.Cells.Item(32, 2).Value = 45923.3125
$ws.Cells.Item(33, 1).Value = 6160
$ws.Cells.Item(33, 2).Value = 45923.32291666666
$ws.Cells.Item(34, 1).Value = 6020
$ws.Cells.Item(34, 2).Value = 45923.33333333334
$ws.Cells.Item(35, 1).Value = 5900
$ws.Cells.Item(35, 2).Value = 45923.34375
$ws.Cells.Item(36, 1).Value = 5770
$ws.Cells.Item(36, 2).Value = 45923.35416666666
$ws.Cells.Item(37, 1).Value = 5630
$ws.Cells.Item(37, 2).Value = 45923.36458333334
$ws.Cells.Item(38, 1).Value = 5510
$ws.Cells.Item(38, 2).Value = 45923.375
$ws.Cells.Item(39, 1).Value = 5370
$ws.Cells.Item(39, 2).Value = 45923.38541666666
$ws.Cells.Item(40, 1).Value = 5250
$ws.Cells.Item(40, 2).Value = 45923.39583333334
$ws.Cells.Item(41, 1).Value = 5140
$ws.Cells.Item(41, 2).Value = 45923.40625
$ws.Cells.Item(42, 1).Value = 5010
$ws.Cells.Item(42, 2).Value = 45923.41666666666
$ws.Cells.Item(43, 1).Value = 4920
$ws.Cells.Item(43, 2).Value = 45923.42708333334
$ws.Cells.Item(44, 1).Value = 4840
$ws.Cells.Item(44, 2).Value = 45923.4375
$ws.Cells.Item(45, 1).Value = 4760
$ws.Cells.Item(45, 2).Value = 45923.44791666666
$ws.Cells.Item(46, 1).Value = 4670
$ws.Cells.Item(46, 2).Value = 45923.45833333334
$ws.Cells.Item(47, 1).Value = 4620
$ws.Cells.Item(47, 2).Value = 45923.46875
$ws.Cells.Item(48, 1).Value = 4600
$ws.Cells.Item(48, 2).Value = 45923.47916666666
$ws.Cells.Item(49, 1).Value = 4600
$ws.Cells.Item(49, 2).Value = 45923.48958333334
$ws.Cells.Item(50, 1).Value = 4610
$ws.Cells.Item(50, 2).Value = 45923.5
$ws.Cells.Item(51, 1).Value = 4620
$ws.Cells.Item(51, 2).Value = 45923.51041666666
$ws.Cells.Item(52, 1).Value = 4630
$ws.Cells.Item(52, 2).Value = 45923.52083333334
$ws.Cells.Item(53, 1).Value = 4650
$ws.Cells.Item(53, 2).Value = 45923.53125
$ws.Cells.Item(54, 1).Value = 4670
$ws.Cells.Item(54, 2).Value = 45923.54166666666
$ws.Cells.Item(55, 1).Value = 4710
$ws.Cells.Item(55, 2).Value = 45923.55208333334
$ws.Cells.Item(56, 1).Value = 4750
$ws.Cells.Item(56, 2).Value = 45923.5625
$ws.Cells.Item(57, 1).Value = 4800
$ws.Cells.Item(57, 2).Value = 45923.57291666666
$ws.Cells.Item(58, 1).Value = 4860
$ws.Cells.Item(58, 2).Value = 45923.58333333334
$ws.Cells.Item(59, 1).Value = 4930
$ws.Cells.Item(59, 2).Value = 45923.59375
$ws.Cells.Item(60, 1).Value = 5000
$ws.Cells.Item(60, 2).Value = 45923.60416666666
$ws.Cells.Item(61, 1).Value = 5090
$ws.Cells.Item(61, 2).Value = 45923.61458333334
$ws.Cells.Item(62, 1).Value = 5190
$ws.Cells.Item(62, 2).Value = 45923.625
$ws.Cells.Item(63, 1).Value = 5310
$ws.Cells.Item(63, 2).Value = 45923.63541666666
$ws.Cells.Item(64, 1).Value = 5440
$ws.Cells.Item(64, 2).Value = 45923.64583333334
$ws.Cells.Item(65, 1).Value = 5590
$ws.Cells.Item(65, 2).Value = 45923.65625
$ws.Cells.Item(66, 1).Value = 5740
$ws.Cells.Item(66, 2).Value = 45923.66666666666
$ws.Cells.Item(67, 1).Value = 5890
$ws.Cells.Item(67, 2).Value = 45923.67708333334
$ws.Cells.Item(68, 1).Value = 6030
$ws.Cells.Item(68, 2).Value = 45923.6875
$ws.Cells.Item(69, 1).Value = 6160
$ws.Cells.Item(69, 2).Value = 45923.69791666666
$ws.Cells.Item(70, 1).Value = 6310
$ws.Cells.Item(70, 2).Value = 45923.70833333334
$ws.Cells.Item(71, 1).Value = 6430
$ws.Cells.Item(71, 2).Value = 45923.71875
$ws.Cells.Item(72, 1).Value = 6560
$ws.Cells.Item(72, 2).Value = 45923.72916666666
$ws.Cells.Item(73, 1).Value = 6680
$ws.Cells.Item(73, 2).Value = 45923.73958333334
$ws.Cells.Item(74, 1).Value = 6800
$ws.Cells.Item(74, 2).Value = 45923.75
$ws.Cells.Item(75, 1).Value = 6930
$ws.Cells.Item(75, 2).Value = 45923.76041666666
$ws.Cells.Item(76, 1).Value = 7040
$ws.Cells.Item(76, 2).Value = 45923.77083333334
$ws.Cells.Item(77, 1).Value = 7130
$ws.Cells.Item(77, 2).Value = 45923.78125
$ws.Cells.Item(78, 1).Value = 7220
$ws.Cells.Item(78, 2).Value = 45923.79166666666
$ws.Cells.Item(79, 1).Value = 7270
$ws.Cells.Item(79, 2).Value = 45923.80208333334
$ws.Cells.Item(80, 1).Value = 7260
$ws.Cells.Item(80, 2).Value = 45923.8125
$ws.Cells.Item(81, 1).Value = 7210
$ws.Cells.Item(81, 2).Value = 45923.82291666666
$ws.Cells.Item(82, 1).Value = 7070
$ws.Cells.Item(82, 2).Value = 45923.83333333334
$ws.Cells.Item(83, 1).Value = 6920
$ws.Cells.Item(83, 2).Value = 45923.84375
$ws.Cells.Item(84, 1).Value = 6790
$ws.Cells.Item(84, 2).Value = 45923.85416666666
$ws.Cells.Item(85, 1).Value = 6670
$ws.Cells.Item(85, 2).Value = 45923.86458333334
$ws.Cells.Item(86, 1).Value = 6520
$ws.Cells.Item(86, 2).Value = 45923.875
$ws.Cells.Item(87, 1).Value = 6390
$ws.Cells.Item(87, 2).Value = 45923.88541666666
$ws.Cells.Item(88, 1).Value = 6230
$ws.Cells.Item(88, 2).Value = 45923.89583333334
$ws.Cells.Item(89, 1).Value = 6060
$ws.Cells.Item(89, 2).Value = 45923.90625
$ws.Cells.Item(90, 1).Value = 5920
$ws.Cells.Item(90, 2).Value = 45923.91666666666
$ws.Cells.Item(91, 1).Value = 5770
$ws.Cells.Item(91, 2).Value = 45923.92708333334
$ws.Cells.Item(92, 1).Value = 5650
$ws.Cells.Item(92, 2).Value = 45923.9375
$ws.Cells.Item(93, 1).Value = 5540
$ws.Cells.Item(93, 2).Value = 45923.94791666666
$ws.Cells.Item(94, 1).Value = 5470
$ws.Cells.Item(94, 2).Value = 45923.95833333334
$ws.Cells.Item(95, 1).Value = 5420
$ws.Cells.Item(95, 2).Value = 45923.96875
$ws.Cells.Item(96, 1).Value = 5370
$ws.Cells.Item(96, 2).Value = 45923.97916666666
$ws.Cells.Item(97, 1).Value = 5310
$ws.Cells.Item(97, 2).Value = 45923.98958333334
